$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-09 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-10 Monday", 2) | Out-Null
$d.Content.Find.Execute("666÷6=111, 0", $true, $false, $false, $false, $false, $true, 1, $false, "454÷6=75, 4", 2) | Out-Null
$d.Content.Find.Execute("543÷6=90, 3", $true, $false, $false, $false, $false, $true, 1, $false, "396÷7=56, 4", 2) | Out-Null
$d.Content.Find.Execute("301÷8=37, 5", $true, $false, $false, $false, $false, $true, 1, $false, "777÷3=259, 0", 2) | Out-Null
$d.Content.Find.Execute("549÷9=61, 0", $true, $false, $false, $false, $false, $true, 1, $false, "923÷2=461, 1", 2) | Out-Null
$d.Content.Find.Execute("506÷7=72, 2", $true, $false, $false, $false, $false, $true, 1, $false, "682÷6=113, 4", 2) | Out-Null
$d.Content.Find.Execute("953÷8=119, 1", $true, $false, $false, $false, $false, $true, 1, $false, "935÷9=103, 8", 2) | Out-Null
$d.Content.Find.Execute("672÷7=96, 0", $true, $false, $false, $false, $false, $true, 1, $false, "277÷5=55, 2", 2) | Out-Null
$d.Content.Find.Execute("962÷5=192, 2", $true, $false, $false, $false, $false, $true, 1, $false, "336÷7=48, 0", 2) | Out-Null
$d.Content.Find.Execute("307÷5=61, 2", $true, $false, $false, $false, $false, $true, 1, $false, "297÷7=42, 3", 2) | Out-Null
$d.Content.Find.Execute("858÷3=286, 0", $true, $false, $false, $false, $false, $true, 1, $false, "726÷9=80, 6", 2) | Out-Null
$d.Content.Find.Execute("591÷8=73, 7", $true, $false, $false, $false, $false, $true, 1, $false, "555÷8=69, 3", 2) | Out-Null
$d.Content.Find.Execute("921÷2=460, 1", $true, $false, $false, $false, $false, $true, 1, $false, "984÷7=140, 4", 2) | Out-Null
$d.Content.Find.Execute("656÷7=93, 5", $true, $false, $false, $false, $false, $true, 1, $false, "833÷3=277, 2", 2) | Out-Null
$d.Content.Find.Execute("711÷5=142, 1", $true, $false, $false, $false, $false, $true, 1, $false, "504÷6=84, 0", 2) | Out-Null
$d.Content.Find.Execute("727÷7=103, 6", $true, $false, $false, $false, $false, $true, 1, $false, "935÷8=116, 7", 2) | Out-Null
$d.Content.Find.Execute("595÷8=74, 3", $true, $false, $false, $false, $false, $true, 1, $false, "488÷3=162, 2", 2) | Out-Null
$d.Content.Find.Execute("174÷9=19, 3", $true, $false, $false, $false, $false, $true, 1, $false, "160÷5=32, 0", 2) | Out-Null
$d.Content.Find.Execute("879÷5=175, 4", $true, $false, $false, $false, $false, $true, 1, $false, "579÷6=96, 3", 2) | Out-Null
$d.Content.Find.Execute("598÷6=99, 4", $true, $false, $false, $false, $false, $true, 1, $false, "161÷9=17, 8", 2) | Out-Null
$d.Content.Find.Execute("645÷2=322, 1", $true, $false, $false, $false, $false, $true, 1, $false, "417÷9=46, 3", 2) | Out-Null
$d.Content.Find.Execute("625÷4=156, 1", $true, $false, $false, $false, $false, $true, 1, $false, "866÷9=96, 2", 2) | Out-Null
$d.Content.Find.Execute("273÷6=45, 3", $true, $false, $false, $false, $false, $true, 1, $false, "650÷4=162, 2", 2) | Out-Null
$d.Content.Find.Execute("578÷8=72, 2", $true, $false, $false, $false, $false, $true, 1, $false, "718÷4=179, 2", 2) | Out-Null
$d.Content.Find.Execute("311÷3=103, 2", $true, $false, $false, $false, $false, $true, 1, $false, "851÷4=212, 3", 2) | Out-Null
$d.Content.Find.Execute("828÷8=103, 4", $true, $false, $false, $false, $false, $true, 1, $false, "807÷7=115, 2", 2) | Out-Null
